# Fill in Group-B trainee rows (8-12) on Sheet1 with names, domains and IDs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names  = @("pankaj", "suraj", "vaibhav", "cyril", "chetan")
$domain = @("unix", "wintel", "wintel", "unix", "automation")
$ids    = @(13309, 13310, 13312, 13313, 13314)

for ($i = 0; $i -lt 5; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 4).Value = $names[$i]   # column D - Name
}

for ($i = 0; $i -lt 5; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 6).Value = $domain[$i]  # column F - Experties
}

for ($i = 0; $i -lt 5; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 7).Value = $ids[$i]     # column G - Group-related id
}

# Excel auto-fits column F's width to the newly entered "Experties" text.
$ws.Columns.Item(6).AutoFit()

# Match the final selection state recorded in the saved workbook.
$ws.Range("F17").Select()
